# Converte as colunas ANO (A) e MES (B) de texto (inlineStr) para numero (t="n"),
# para as linhas de dados 2..117. O ANO permanece 2021 (numero) e o MES perde o
# zero a esquerda (ex.: "02" -> 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 117

for ($r = 2; $r -le $lastRow; $r++) {
    $anoCell = $ws.Cells.Item($r, 1)
    $mesCell = $ws.Cells.Item($r, 2)

    $anoNum = $anoCell.Text + 0
    $mesNum = $mesCell.Text + 0

    $anoCell.Value = $anoNum
    $mesCell.Value = $mesNum
}
